$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.426.85'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.941.99'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.47'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.610'
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.13'
$ws.Range('E8').Value = '  -2.72%  '
$ws.Range('E9').Value = '  -3.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0850'
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('D12').Value = '2.226.20'
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.810'
$ws.Range('E13').Value = '  -5.10%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.18'
$ws.Range('E14').Value = '  -3.53%  '
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('E16').Value = '  -5.65%  '
$ws.Range('D17').Value = '1.942.74'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('D18').Value = '36.399.21'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.16'
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('D20').Value = '0.0₃0863'
$ws.Range('E20').Value = '  -5.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '227.81'
$ws.Range('E22').Value = '  -5.00%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.34'
$ws.Range('E24').Value = '  -6.58%  '
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.20'
$ws.Range('E26').Value = '  -5.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.60'
$ws.Range('E27').Value = '  -1.99%  '
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.22'
$ws.Range('E29').Value = '  -3.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.118'
$ws.Range('E30').Value = '  -1.61%  '
$ws.Range('E31').Value = '  -6.95%  '
$ws.Range('E32').Value = '  -6.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0617'
$ws.Range('E33').Value = '  -4.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.16'
$ws.Range('E34').Value = '  -5.55%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.15'
$ws.Range('E39').Value = '  +8.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0987'
$ws.Range('E40').Value = '  +2.45%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0210'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('E43').Value = '  -5.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.73'
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('D45').Value = '1.339.65'
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('E46').Value = '  -5.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.42'
$ws.Range('E47').Value = '  -4.94%  '
$ws.Range('E48').Value = '  -3.74%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = '2.116.68'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.99'
$ws.Range('E51').Value = '  -4.96%  '
